# Apply the "add excel / json field" change: a new boolean-ish "sold" column
# (header + true/false text values) appended as column J on the "Json" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Json")

# --- New column header ------------------------------------------------
$ws.Range("J1").Value = "sold"

# --- New column values --------------------------------------------------
# Plain "true"/"false" need to land in the sheet as literal text (shared
# string), not as native Excel booleans. Writing the text through a
# formula and then collapsing the formula to its value via copy/paste
# values keeps the result a plain inline string cell with no extra
# number-format/style baggage.
$ws.Range("J2").Formula = '="true"'
$ws.Range("J2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

$ws.Range("J3").Formula = '="false"'
$ws.Range("J3").Copy()
$ws.Range("J3").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- Refresh the view / selection ---------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("J4").Select()
